$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (in-place shared-string edits) ---
$ws.Cells.Replace("20.01.2025", "28.01.2025")
$ws.Cells.Replace("(20/01/2025)", "(29/01/2025)")

# --- Top block (rows 3-6) ---
$ws.Range("C3").Value = 43137
$ws.Range("D3").Value = 7623

$ws.Range("C4").Value = 29436
$ws.Range("D4").Value = 9975
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()

$ws.Range("C5").Value = 38971
$ws.Range("D5").Value = 17598

$ws.Range("C6").Value = 36591
$ws.Range("D6").Value = 6964
$ws.Range("E6").ClearContents()

# --- Lower block ---
$ws.Range("C13").Value = 104517

$ws.Range("C14").Value = 272209
$ws.Range("D14").Value = 148135
$ws.Range("E14").Value = 114285

$ws.Range("C20").Value = 5850
$ws.Range("D20").Value = 2060

$ws.Range("C21").Value = 310
$ws.Range("D21").Value = 120

$ws.Range("C22").Value = 590
$ws.Range("D22").ClearContents()

$ws.Range("C24").Value = 23
$ws.Range("D24").ClearContents()

$ws.Range("C26").Value = 44
$ws.Range("D26").ClearContents()

$ws.Range("C27").Value = 84
$ws.Range("D27").Value = 1

# --- Selection (view state) ---
$ws.Range("J30").Select()
